$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("I2").Value = 0.7635196712427992
$ws.Range("J2").Value = 0.7635196712427992
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 6.297819150816
$ws.Range("R2").Value = 56.68037235734399
$ws.Range("S2").Value = 0.5234401701313929
$ws.Range("T2").Value = 0.5234401701313929

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.7635196712427992
$ws.Range("J3").Value = 0.7635196712427992
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("S3").Value = 0.1808217742350918
$ws.Range("T3").Value = 0.1808217742350918

# Row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 0.7635196712427992
$ws.Range("J4").Value = 0.7635196712427992
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 0.7129648591200001
$ws.Range("R4").Value = 6.41668373208
$ws.Range("S4").Value = 0.05925772687631449
$ws.Range("T4").Value = 0.0592577268763145

# Row 5 (MuSCs -> ECs)
$ws.Range("G5").Value = 0.519749
$ws.Range("H5").Value = 1.559247
$ws.Range("I5").Value = 0.2364803287572008
$ws.Range("J5").Value = 0.2364803287572008
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.752937333333333
$ws.Range("N5").Value = 11.258812
$ws.Range("O5").Value = 0.6855621274031838
$ws.Range("P5").Value = 0.6855621274031838
$ws.Range("Q5").Value = 1.950585426062667
$ws.Range("R5").Value = 17.555268834564
$ws.Range("S5").Value = 0.1621219572717909
$ws.Range("T5").Value = 0.1621219572717909

# Row 6 (MuSCs -> FAPs)
$ws.Range("G6").Value = 0.519749
$ws.Range("H6").Value = 1.559247
$ws.Range("I6").Value = 0.2364803287572008
$ws.Range("J6").Value = 0.2364803287572008
$ws.Range("O6").Value = 0.2368266084628361
$ws.Range("P6").Value = 0.2368266084628362
$ws.Range("Q6").Value = 0.6738273783023333
$ws.Range("R6").Value = 6.064446404721001
$ws.Range("S6").Value = 0.05600483422774435
$ws.Range("T6").Value = 0.05600483422774437

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = 0.519749
$ws.Range("H7").Value = 1.559247
$ws.Range("I7").Value = 0.2364803287572008
$ws.Range("J7").Value = 0.2364803287572008
$ws.Range("M7").Value = 0.4248633333333334
$ws.Range("N7").Value = 1.27459
$ws.Range("O7").Value = 0.07761126413398003
$ws.Range("P7").Value = 0.07761126413398005
$ws.Range("Q7").Value = 0.2208222926366667
$ws.Range("R7").Value = 1.98740063373
$ws.Range("S7").Value = 0.01835353725766554
$ws.Range("T7").Value = 0.01835353725766555
